$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "space jump"
$ws.Range("B8").Value = "com.singleton.helix"

$ws.Range("B13").Select()
